$wb = $excel.ActiveWorkbook

# --- Update text on Sheet1-Merged (rows 2-4): drop "--newSheet=Sheet2", drop "./" prefix, fix "Apend" -> "Append" ---
$ws1 = $wb.Worksheets.Item("Sheet1-Merged")
$ws1.Range("A2").Value = "Diff test:  xltablediff.py  --key=ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"
$ws1.Range("A3").Value = "Merge test:  xltablediff.py  --key=ID test1old.xlsx test1new.xlsx --merge=Color --out test1merge.xlsx"
$ws1.Range("A4").Value = "Append test:  xltablediff.py  --key=ID test1old.xlsx test1new.xlsx --append --out test1append.xlsx"

# --- Clear the trailing row 19 contents (was "Trailing row here" in A19 and a numeric cell in F19) ---
$ws1.Range("A19:F19").Clear()

# --- Update selection / active cell on Sheet1-Merged ---
$ws1.Range("F17").Select() | Out-Null

# --- Delete Sheet2 ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete() | Out-Null

# --- Make Sheet1-Merged the active tab ---
$ws1.Activate() | Out-Null
